$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "33×63=" "87×37="
Replace-Text "53×28=" "27×32="
Replace-Text "63×54=" "42×84="
Replace-Text "36×74=" "58×86="
Replace-Text "25×79=" "61×86="
Replace-Text "28×19=" "22×42="
Replace-Text "71×81=" "69×54="
Replace-Text "95×55=" "33×17="
Replace-Text "89×58=" "26×85="
Replace-Text "18×28=" "87×33="
Replace-Text "39×40=" "52×23="
Replace-Text "64×75=" "52×57="
Replace-Text "83×72=" "56×71="
Replace-Text "51×47=" "64×97="
Replace-Text "17×93=" "69×28="
Replace-Text "61×77=" "68×70="
Replace-Text "29×97=" "67×68="
Replace-Text "85×48=" "11×33="
Replace-Text "47×96=" "98×45="
Replace-Text "48×92=" "77×15="
Replace-Text "52×84=" "14×92="
Replace-Text "90×97=" "60×97="
Replace-Text "70×39=" "34×50="
Replace-Text "54×23=" "32×31="
Replace-Text "85×71=" "21×47="
